$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) values for rows 2-11

$ws.Range("C2").Value = -0.1679909953156797
$ws.Range("D2").Value = 0.8675856690834247

$ws.Range("C3").Value = -0.8927145915148792
$ws.Range("D3").Value = 0.3782852011271101

$ws.Range("C4").Value = -2.233523161090723
$ws.Range("D4").Value = 0.03220112947831488

$ws.Range("C5").Value = -2.964546938285084
$ws.Range("D5").Value = 0.005506704071783775

$ws.Range("C6").Value = -0.5813551461600223
$ws.Range("D6").Value = 0.5648356671766237

$ws.Range("C7").Value = -2.150838376735805
$ws.Range("D7").Value = 0.03868434121037923

$ws.Range("C8").Value = -2.650015155701639
$ws.Range("D8").Value = 0.01212704975196344

$ws.Range("C9").Value = -2.023892997956647
$ws.Range("D9").Value = 0.05089661164679593

$ws.Range("C10").Value = -2.40616272574227
$ws.Range("D10").Value = 0.02170422856365617

$ws.Range("C11").Value = -0.6458905805875922
$ws.Range("D11").Value = 0.5226859634931962

$wb.Save()
